$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("nomor induk") shifting username/nama/password right
$ws.Columns("B").Insert()

# Format the new column as Text so the long ID numbers are kept verbatim
$ws.Columns("B").NumberFormat = "@"
$ws.Columns("B").ColumnWidth = 23.5

# Re-assert the original column widths for the (now shifted) columns
$ws.Columns("C").ColumnWidth = 14.5
$ws.Columns("D").ColumnWidth = 10.333333333333334

# Header
$ws.Range("B1").Value = "nomor induk"

# Data - add B3 before B2 so shared-string insertion order matches the source workbook
$ws.Range("B3").Value = "198311052003101000"
$ws.Range("B2").Value = "195912312010123890"

# Match the saved selection from the source workbook
$ws.Range("I17").Select()
